# Weekly update: a new "Haba" price record for
# "Terminal Hortofrutícola Agro Chillán" is inserted as row 10, pushing the
# existing rows 10-22 down to rows 11-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 (shifts rows 10..22 -> 11..23,
# and extends the used range from A1:R22 to A1:R23).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's data point.
$ws.Range('A10').Value = 7
$ws.Range('B10').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C10').Value = 'Ñuble'
$ws.Range('D10').Value = 44512
$ws.Range('E10').Value = 16
$ws.Range('F10').Value = 100112026
$ws.Range('G10').Value = 'Haba'
$ws.Range('H10').Value = 'Sin especificar'
$ws.Range('I10').Value = 'Primera'
$ws.Range('J10').Value = 100
$ws.Range('K10').Value = 7000
$ws.Range('L10').Value = 8000
$ws.Range('M10').Value = 7500
$ws.Range('N10').Value = '$/saco 25 kilos'
$ws.Range('O10').Value = 'Provincia de Diguillín'
$ws.Range('P10').Value = 300
$ws.Range('Q10').Value = 25
$ws.Range('R10').Value = 'Hortaliza'
